$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (non-numeric) storage for the Price/Volume columns so that
# numeric- and percent-looking strings are not auto-converted to numbers.
$dataRange = $ws.Range("D2:E47")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "274.58"
$ws.Range("E2").Value = "0.42%"
$ws.Range("E3").Value = "2.06%"
$ws.Range("D4").Value = "4.878"
$ws.Range("E4").Value = "0.12%"
$ws.Range("D5").Value = "0.06391"
$ws.Range("D6").Value = "6.939"
$ws.Range("E6").Value = "0.92%"
$ws.Range("D7").Value = "1.188"
$ws.Range("E7").Value = "-1.11%"
$ws.Range("D8").Value = "0.8769"
$ws.Range("E8").Value = "0.74%"
$ws.Range("D9").Value = "0.1513"
$ws.Range("E9").Value = "3.81%"
$ws.Range("D10").Value = "0.05110"
$ws.Range("E10").Value = "0.02%"
$ws.Range("D11").Value = "0.07534"
$ws.Range("E11").Value = "2.14%"
$ws.Range("D12").Value = "0.02964"
$ws.Range("E12").Value = "-1.54%"
$ws.Range("D13").Value = "0.08974"
$ws.Range("E13").Value = "-0.75%"
$ws.Range("D14").Value = "0.001562"
$ws.Range("E14").Value = "-0.84%"
$ws.Range("D15").Value = "0.0006381"
$ws.Range("E15").Value = "1.00%"
$ws.Range("D16").Value = "0.006185"
$ws.Range("E16").Value = "2.38%"
$ws.Range("D17").Value = "3.467"
$ws.Range("E17").Value = "0.38%"
$ws.Range("D18").Value = "3.313"
$ws.Range("E18").Value = "-1.01%"
$ws.Range("E19").Value = "0.51%"
$ws.Range("E21").Value = "1.80%"
$ws.Range("D22").Value = "3.923"
$ws.Range("E22").Value = "-0.36%"
$ws.Range("D23").Value = "0.04411"
$ws.Range("E23").Value = "0.21%"
$ws.Range("D25").Value = "0.001178"
$ws.Range("E25").Value = "0.14%"
$ws.Range("D26").Value = "0.003853"
$ws.Range("E26").Value = "-9.63%"
$ws.Range("E27").Value = "0.08%"
$ws.Range("E28").Value = "14.57%"
$ws.Range("E40").Value = "2.42%"
$ws.Range("D41").Value = "0.006798"
$ws.Range("E41").Value = "1.19%"
$ws.Range("E42").Value = "0.79%"
$ws.Range("D43").Value = "0.002192"
$ws.Range("E43").Value = "4.37%"
$ws.Range("D44").Value = "0.01148"
$ws.Range("E44").Value = "-8.34%"
$ws.Range("D45").Value = "0.00005197"
$ws.Range("E45").Value = "-1.90%"
$ws.Range("D46").Value = "1.681"
$ws.Range("E46").Value = "-43.75%"
$ws.Range("E47").Value = "0.02%"

# Restore the default cell style so no stray number-format style lingers
# on the edited cells (matches original formatting of the workbook).
$dataRange.Style = "Normal"

